$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2"  "267.86"
Set-TextValue "D3"  "22.87"
Set-TextValue "D5"  "0.06202"
Set-TextValue "D7"  "6.672"
Set-TextValue "D8"  "1.391"
Set-TextValue "D9"  "0.8295"
Set-TextValue "D11" "0.1609"
Set-TextValue "D12" "0.08229"
Set-TextValue "D13" "0.03395"
Set-TextValue "D14" "0.03153"
Set-TextValue "D15" "0.09286"
Set-TextValue "D16" "3.921"
Set-TextValue "D17" "0.001717"
Set-TextValue "D18" "0.04848"
Set-TextValue "D19" "0.006298"
Set-TextValue "D20" "0.005387"
Set-TextValue "D21" "0.001091"
Set-TextValue "D23" "3.757"
Set-TextValue "D24" "2.367"
Set-TextValue "D26" "0.1212"
Set-TextValue "D40" "0.04656"
Set-TextValue "D41" "0.006880"
Set-TextValue "D42" "0.1155"
Set-TextValue "D43" "0.003350"
Set-TextValue "D44" "0.01225"
Set-TextValue "D45" "0.00006265"
Set-TextValue "D48" "0.1637"
